$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -20.02482043281714
$ws.Cells.Item(2, 3).Value = 1.958990465535035
$ws.Cells.Item(2, 4).Value = -20.02482043281714
$ws.Cells.Item(2, 5).Value = -20.02482043281714
$ws.Cells.Item(2, 6).Value = -20.02482043281714
$ws.Cells.Item(2, 7).Value = -20.02482043281714
$ws.Cells.Item(2, 8).Value = -20.02482043281714
$ws.Cells.Item(2, 9).Value = -20.02482043281714
$ws.Cells.Item(2, 10).Value = -20.02482043281714
$ws.Cells.Item(2, 11).Value = -20.02482043281714

$ws.Cells.Item(3, 2).Value = -20.02482043281714
$ws.Cells.Item(3, 3).Value = -20.02482043281714
$ws.Cells.Item(3, 4).Value = -20.02482043281714
$ws.Cells.Item(3, 5).Value = -20.02482043281714
$ws.Cells.Item(3, 6).Value = -20.02482043281714
$ws.Cells.Item(3, 7).Value = -20.02482043281714
$ws.Cells.Item(3, 8).Value = -20.02482043281714
$ws.Cells.Item(3, 9).Value = 1.430701950586065
$ws.Cells.Item(3, 10).Value = -20.02482043281714
$ws.Cells.Item(3, 11).Value = -20.02482043281714

$ws.Cells.Item(4, 2).Value = -20.02482043281714
$ws.Cells.Item(4, 3).Value = 1.96812377943874
$ws.Cells.Item(4, 4).Value = 1.667625226896033
$ws.Cells.Item(4, 5).Value = -20.02482043281714
$ws.Cells.Item(4, 6).Value = 3.429444341277538
$ws.Cells.Item(4, 7).Value = -20.02482043281714
$ws.Cells.Item(4, 8).Value = 1.28463871235626
$ws.Cells.Item(4, 9).Value = -20.02482043281714
$ws.Cells.Item(4, 10).Value = 0.7433194597944037
$ws.Cells.Item(4, 11).Value = -20.02482043281714

$ws.Cells.Item(5, 2).Value = -20.02482043281714
$ws.Cells.Item(5, 3).Value = 1.691378166530212
$ws.Cells.Item(5, 4).Value = -20.02482043281714
$ws.Cells.Item(5, 5).Value = -20.02482043281714
$ws.Cells.Item(5, 6).Value = -20.02482043281714
$ws.Cells.Item(5, 7).Value = -20.02482043281714
$ws.Cells.Item(5, 8).Value = -20.02482043281714
$ws.Cells.Item(5, 9).Value = -20.02482043281714
$ws.Cells.Item(5, 10).Value = -20.02482043281714
$ws.Cells.Item(5, 11).Value = -20.02482043281714

$ws.Cells.Item(6, 2).Value = -20.02482043281714
$ws.Cells.Item(6, 3).Value = -20.02482043281714
$ws.Cells.Item(6, 4).Value = -20.02482043281714
$ws.Cells.Item(6, 5).Value = -20.02482043281714
$ws.Cells.Item(6, 6).Value = -20.02482043281714
$ws.Cells.Item(6, 7).Value = -20.02482043281714
$ws.Cells.Item(6, 8).Value = -20.02482043281714
$ws.Cells.Item(6, 9).Value = -20.02482043281714
$ws.Cells.Item(6, 10).Value = -20.02482043281714
$ws.Cells.Item(6, 11).Value = -20.02482043281714

$ws.Cells.Item(7, 2).Value = 2.466137403366719
$ws.Cells.Item(7, 3).Value = -20.02482043281714
$ws.Cells.Item(7, 4).Value = -20.02482043281714
$ws.Cells.Item(7, 5).Value = -20.02482043281714
$ws.Cells.Item(7, 6).Value = -20.02482043281714
$ws.Cells.Item(7, 7).Value = -20.02482043281714
$ws.Cells.Item(7, 8).Value = -20.02482043281714
$ws.Cells.Item(7, 9).Value = -20.02482043281714
$ws.Cells.Item(7, 10).Value = -20.02482043281714
$ws.Cells.Item(7, 11).Value = -20.02482043281714

$ws.Cells.Item(8, 2).Value = -20.02482043281714
$ws.Cells.Item(8, 3).Value = -20.02482043281714
$ws.Cells.Item(8, 4).Value = -20.02482043281714
$ws.Cells.Item(8, 5).Value = 1.778032336596316
$ws.Cells.Item(8, 6).Value = -20.02482043281714
$ws.Cells.Item(8, 7).Value = -20.02482043281714
$ws.Cells.Item(8, 8).Value = -20.02482043281714
$ws.Cells.Item(8, 9).Value = -20.02482043281714
$ws.Cells.Item(8, 10).Value = -20.02482043281714
$ws.Cells.Item(8, 11).Value = -20.02482043281714

$ws.Cells.Item(9, 2).Value = 3.855427563889433
$ws.Cells.Item(9, 3).Value = -20.02482043281714
$ws.Cells.Item(9, 4).Value = -20.02482043281714
$ws.Cells.Item(9, 5).Value = -20.02482043281714
$ws.Cells.Item(9, 6).Value = -20.02482043281714
$ws.Cells.Item(9, 7).Value = -20.02482043281714
$ws.Cells.Item(9, 8).Value = -20.02482043281714
$ws.Cells.Item(9, 9).Value = -20.02482043281714
$ws.Cells.Item(9, 10).Value = -20.02482043281714
$ws.Cells.Item(9, 11).Value = -20.02482043281714

$ws.Cells.Item(10, 2).Value = -20.02482043281714
$ws.Cells.Item(10, 3).Value = -20.02482043281714
$ws.Cells.Item(10, 4).Value = -20.02482043281714
$ws.Cells.Item(10, 5).Value = -20.02482043281714
$ws.Cells.Item(10, 6).Value = -20.02482043281714
$ws.Cells.Item(10, 7).Value = -20.02482043281714
$ws.Cells.Item(10, 8).Value = -20.02482043281714
$ws.Cells.Item(10, 9).Value = 1.642833851248302
$ws.Cells.Item(10, 10).Value = -20.02482043281714
$ws.Cells.Item(10, 11).Value = 2.315253767065207

$ws.Cells.Item(11, 2).Value = -20.02482043281714
$ws.Cells.Item(11, 3).Value = -20.02482043281714
$ws.Cells.Item(11, 4).Value = -20.02482043281714
$ws.Cells.Item(11, 5).Value = 2.930347884484047
$ws.Cells.Item(11, 6).Value = -20.02482043281714
$ws.Cells.Item(11, 7).Value = 4.321926810113391
$ws.Cells.Item(11, 8).Value = -20.02482043281714
$ws.Cells.Item(11, 9).Value = -20.02482043281714
$ws.Cells.Item(11, 10).Value = -20.02482043281714
$ws.Cells.Item(11, 11).Value = 2.052934439121146

$ws.Cells.Item(12, 2).Value = -20.02482043281714
$ws.Cells.Item(12, 3).Value = -20.02482043281714
$ws.Cells.Item(12, 4).Value = -20.02482043281714
$ws.Cells.Item(12, 5).Value = -20.02482043281714
$ws.Cells.Item(12, 6).Value = -20.02482043281714
$ws.Cells.Item(12, 7).Value = -20.02482043281714
$ws.Cells.Item(12, 8).Value = -20.02482043281714
$ws.Cells.Item(12, 9).Value = -20.02482043281714
$ws.Cells.Item(12, 10).Value = -20.02482043281714
$ws.Cells.Item(12, 11).Value = -20.02482043281714

$ws.Cells.Item(13, 2).Value = -20.02482043281714
$ws.Cells.Item(13, 3).Value = -20.02482043281714
$ws.Cells.Item(13, 4).Value = -20.02482043281714
$ws.Cells.Item(13, 5).Value = 2.521853872480502
$ws.Cells.Item(13, 6).Value = -20.02482043281714
$ws.Cells.Item(13, 7).Value = -20.02482043281714
$ws.Cells.Item(13, 8).Value = -20.02482043281714
$ws.Cells.Item(13, 9).Value = -20.02482043281714
$ws.Cells.Item(13, 10).Value = 1.755154149394715
$ws.Cells.Item(13, 11).Value = 1.777716836614956

$ws.Cells.Item(14, 2).Value = -20.02482043281714
$ws.Cells.Item(14, 3).Value = -20.02482043281714
$ws.Cells.Item(14, 4).Value = 1.526246133308786
$ws.Cells.Item(14, 5).Value = -20.02482043281714
$ws.Cells.Item(14, 6).Value = -20.02482043281714
$ws.Cells.Item(14, 7).Value = -20.02482043281714
$ws.Cells.Item(14, 8).Value = -20.02482043281714
$ws.Cells.Item(14, 9).Value = -20.02482043281714
$ws.Cells.Item(14, 10).Value = -20.02482043281714
$ws.Cells.Item(14, 11).Value = 1.880528741977888

$ws.Cells.Item(15, 2).Value = -20.02482043281714
$ws.Cells.Item(15, 3).Value = -20.02482043281714
$ws.Cells.Item(15, 4).Value = 1.773964990377483
$ws.Cells.Item(15, 5).Value = -20.02482043281714
$ws.Cells.Item(15, 6).Value = -20.02482043281714
$ws.Cells.Item(15, 7).Value = -20.02482043281714
$ws.Cells.Item(15, 8).Value = -20.02482043281714
$ws.Cells.Item(15, 9).Value = -20.02482043281714
$ws.Cells.Item(15, 10).Value = -20.02482043281714
$ws.Cells.Item(15, 11).Value = -20.02482043281714

$ws.Cells.Item(16, 2).Value = -20.02482043281714
$ws.Cells.Item(16, 3).Value = -20.02482043281714
$ws.Cells.Item(16, 4).Value = -20.02482043281714
$ws.Cells.Item(16, 5).Value = -20.02482043281714
$ws.Cells.Item(16, 6).Value = -20.02482043281714
$ws.Cells.Item(16, 7).Value = -20.02482043281714
$ws.Cells.Item(16, 8).Value = -20.02482043281714
$ws.Cells.Item(16, 9).Value = -20.02482043281714
$ws.Cells.Item(16, 10).Value = 1.888959045016571
$ws.Cells.Item(16, 11).Value = -20.02482043281714

$ws.Cells.Item(17, 2).Value = -20.02482043281714
$ws.Cells.Item(17, 3).Value = 2.134471039941486
$ws.Cells.Item(17, 4).Value = 1.833446333588669
$ws.Cells.Item(17, 5).Value = -20.02482043281714
$ws.Cells.Item(17, 6).Value = -20.02482043281714
$ws.Cells.Item(17, 7).Value = -20.02482043281714
$ws.Cells.Item(17, 8).Value = 2.084570596090808
$ws.Cells.Item(17, 9).Value = 2.204309645669665
$ws.Cells.Item(17, 10).Value = 2.558431546721403
$ws.Cells.Item(17, 11).Value = -20.02482043281714

$ws.Cells.Item(18, 2).Value = -20.02482043281714
$ws.Cells.Item(18, 3).Value = -20.02482043281714
$ws.Cells.Item(18, 4).Value = -20.02482043281714
$ws.Cells.Item(18, 5).Value = -20.02482043281714
$ws.Cells.Item(18, 6).Value = -20.02482043281714
$ws.Cells.Item(18, 7).Value = -20.02482043281714
$ws.Cells.Item(18, 8).Value = 2.056437059364739
$ws.Cells.Item(18, 9).Value = 2.016349737911555
$ws.Cells.Item(18, 10).Value = 2.421153478074485
$ws.Cells.Item(18, 11).Value = -20.02482043281714

$ws.Cells.Item(19, 2).Value = -20.02482043281714
$ws.Cells.Item(19, 3).Value = -20.02482043281714
$ws.Cells.Item(19, 4).Value = 2.037811082641815
$ws.Cells.Item(19, 5).Value = -20.02482043281714
$ws.Cells.Item(19, 6).Value = -20.02482043281714
$ws.Cells.Item(19, 7).Value = -20.02482043281714
$ws.Cells.Item(19, 8).Value = 1.541025882551333
$ws.Cells.Item(19, 9).Value = 1.749791712181
$ws.Cells.Item(19, 10).Value = -20.02482043281714
$ws.Cells.Item(19, 11).Value = -20.02482043281714

$ws.Cells.Item(20, 2).Value = -20.02482043281714
$ws.Cells.Item(20, 3).Value = 1.052778316045752
$ws.Cells.Item(20, 4).Value = 1.512328280427457
$ws.Cells.Item(20, 5).Value = -20.02482043281714
$ws.Cells.Item(20, 6).Value = 3.205746706171197
$ws.Cells.Item(20, 7).Value = -20.02482043281714
$ws.Cells.Item(20, 8).Value = 1.766733979538702
$ws.Cells.Item(20, 9).Value = 1.113829478198912
$ws.Cells.Item(20, 10).Value = -20.02482043281714
$ws.Cells.Item(20, 11).Value = 1.911745982092775

$ws.Cells.Item(21, 2).Value = -20.02482043281714
$ws.Cells.Item(21, 3).Value = 1.324560694092738
$ws.Cells.Item(21, 4).Value = -20.02482043281714
$ws.Cells.Item(21, 5).Value = 1.679993811840832
$ws.Cells.Item(21, 6).Value = -20.02482043281714
$ws.Cells.Item(21, 7).Value = -20.02482043281714
$ws.Cells.Item(21, 8).Value = 1.510895864041396
$ws.Cells.Item(21, 9).Value = -20.02482043281714
$ws.Cells.Item(21, 10).Value = -20.02482043281714
$ws.Cells.Item(21, 11).Value = -20.02482043281714

